# "Drop in results from RMI script"
#
#   - removes the "Texas Data" worksheet entirely (its cells were the
#     only place the explanatory notes about the NREL waste-heat
#     calculation error were used, so those shared strings disappear
#     along with the sheet -- no separate edit needed for that)
#   - restores the HPEbP natural-gas-reforming efficiency formula in
#     B3 to include waste heat in the energy-balance denominator
#     (118 / (162 + 2 + 46) instead of 118 / (162 + 2)); every other
#     cell in that row (C3:AI3) is a copied/shared formula that just
#     references back to B3, so it follows along automatically
#   - restores the previous active sheet / selections on each tab

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove the "Texas Data" worksheet entirely -----------------------
$wsTexas = $wb.Worksheets.Item("Texas Data")
$wsTexas.Delete()

# --- HPEbP: restore the old (waste-heat-inclusive) efficiency formula -
$wsHPEbP = $wb.Worksheets.Item("HPEbP")
$wsHPEbP.Range("B3").Formula = "=118/(162+2+46)"

# --- Restore prior selections / active sheet ---------------------------
$wsIEA = $wb.Worksheets.Item("IEA Data")
$wsIEA.Activate()
$wsIEA.Range("A30").Select()

$wsHPEbP.Activate()
$wsHPEbP.Range("H17:H18").Select()

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("B16").Select()
